$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Write cell values in reading order (A,B,C,D per row, row by row)
#    so new shared-string entries land in the same order as the target.
# ---------------------------------------------------------------
$ws.Range('A19').Value = 'Entity1'
$ws.Range('B19').Value = 'Entity2'
$ws.Range('C19').Value = 'Relation'
$ws.Range('D19').Value = 'Score'
$ws.Range('A20').Value = 'sherlock holmes'
$ws.Range('B20').Value = 'united kingdom'
$ws.Range('C20').Value = 'anthem'
$ws.Range('D20').Value = [double]'7.0793549999999998E-3'
$ws.Range('C21').Value = 'firstAppearance'
$ws.Range('D21').Value = [double]'4.4867869999999999E-3'
$ws.Range('C22').Value = 'allegiance'
$ws.Range('D22').Value = [double]'4.3567459999999999E-3'
$ws.Range('A23').Value = 'apple inc'
$ws.Range('B23').Value = 'steve jobs'
$ws.Range('C23').Value = 'foundedBy'
$ws.Range('D23').Value = [double]'1.5438646E-2'
$ws.Range('C24').Value = 'keyPerson'
$ws.Range('D24').Value = [double]'9.9318410000000003E-3'
$ws.Range('C25').Value = 'successor'
$ws.Range('D25').Value = [double]'8.0690020000000005E-3'
$ws.Range('A26').Value = 'adolf hitler'
$ws.Range('B26').Value = 'world war ii'
$ws.Range('C26').Value = 'commander'
$ws.Range('D26').Value = [double]'3.7711715999999999E-2'
$ws.Range('C27').Value = 'battle'
$ws.Range('D27').Value = [double]'2.2160928999999999E-2'
$ws.Range('C28').Value = 'ceo'
$ws.Range('D28').Value = [double]'2.44E-5'
$ws.Range('A29').Value = 'microsoft'
$ws.Range('B29').Value = 'redmond'
$ws.Range('C29').Value = 'locationCity'
$ws.Range('D29').Value = [double]'8.2507291999999996E-2'
$ws.Range('C30').Value = 'foundationPlace'
$ws.Range('D30').Value = [double]'4.7192255000000002E-2'
$ws.Range('C31').Value = 'location'
$ws.Range('D31').Value = [double]'3.6915687000000003E-2'
$ws.Range('A32').Value = 'titanic'
$ws.Range('B32').Value = 'james cameron'
$ws.Range('C32').Value = 'director'
$ws.Range('D32').Value = [double]'0.12440701799999999'
$ws.Range('C33').Value = 'cinematography'
$ws.Range('D33').Value = [double]'9.6447235000000006E-2'
$ws.Range('C34').Value = 'editing'
$ws.Range('D34').Value = [double]'8.0133909000000003E-2'
$ws.Range('A35').Value = 'titanic'
$ws.Range('B35').Value = 'leonardo dicaprio'
$ws.Range('C35').Value = 'starring'
$ws.Range('D35').Value = [double]'4.9688828999999997E-2'
$ws.Range('C36').Value = 'narrator'
$ws.Range('D36').Value = [double]'3.7266747000000003E-2'
$ws.Range('C37').Value = 'producer'
$ws.Range('D37').Value = [double]'1.3059708E-2'
$ws.Range('A38').Value = 'harry potter'
$ws.Range('B38').Value = 'j k rowling'
$ws.Range('C38').Value = 'notableWork'
$ws.Range('D38').Value = [double]'1.6964505000000001E-2'
$ws.Range('C39').Value = 'author'
$ws.Range('D39').Value = [double]'1.5514061000000001E-2'
$ws.Range('C40').Value = 'coverArtist'
$ws.Range('D40').Value = [double]'1.4906320000000001E-2'

# ---------------------------------------------------------------
# 2) Formatting - apply in the same order the target style table was
#    built so the generated cellXfs land at indexes 3, 4, 5.
# ---------------------------------------------------------------

# style index 3: wrap text only (header row + all Relation/Score cells)
$ws.Range('A19:D19').WrapText = $true

# style index 4: wrap text + centered (Entity1/Entity2 columns, rows 20-40)
$ws.Range('A20:B40').WrapText = $true
$ws.Range('A20:B40').HorizontalAlignment = -4108

# style index 3 again: Relation/Score columns, rows 20-40
$ws.Range('C20:D40').WrapText = $true

# style index 5: D28 (ceo score) uses scientific notation + wrap
$ws.Range('D28').WrapText = $true
$ws.Range('D28').NumberFormat = "0.00E+00"

# ---------------------------------------------------------------
# 3) Merge the Entity1/Entity2 cells that span multiple relation rows
# ---------------------------------------------------------------
$ws.Range('A20:A22').Merge()
$ws.Range('B20:B22').Merge()
$ws.Range('A23:A25').Merge()
$ws.Range('B23:B25').Merge()
$ws.Range('A26:A28').Merge()
$ws.Range('B26:B28').Merge()
$ws.Range('A29:A31').Merge()
$ws.Range('B29:B31').Merge()
$ws.Range('A32:A34').Merge()
$ws.Range('B32:B34').Merge()
$ws.Range('A35:A37').Merge()
$ws.Range('B35:B37').Merge()
$ws.Range('A38:A40').Merge()
$ws.Range('B38:B40').Merge()

# ---------------------------------------------------------------
# 4) Row heights for wrapped two-line rows
# ---------------------------------------------------------------
$ws.Rows.Item(21).RowHeight = 27
$ws.Rows.Item(22).RowHeight = 27
$ws.Rows.Item(23).RowHeight = 27
$ws.Rows.Item(24).RowHeight = 27
$ws.Rows.Item(25).RowHeight = 27
$ws.Rows.Item(26).RowHeight = 27
$ws.Rows.Item(29).RowHeight = 27
$ws.Rows.Item(30).RowHeight = 27
$ws.Rows.Item(33).RowHeight = 27
$ws.Rows.Item(38).RowHeight = 27
$ws.Rows.Item(40).RowHeight = 27

# ---------------------------------------------------------------
# 5) Selection matches the authored state after adding the table
# ---------------------------------------------------------------
$ws.Range('A19:D40').Select()
